$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Zagłębie Lubin"
$ws.Range("A11").Value = "Remis"
$ws.Range("A15").Value = "Zagłębie Lubin"
$ws.Range("A17").Value = "Radomiak Radom"
$ws.Range("A18").Value = "Remis"
$ws.Range("A20").Value = "Legia Warszawa"
$ws.Range("A26").Value = "Górnik Zabrze"
$ws.Range("A27").Value = "Remis"
$ws.Range("A29").Value = "Jagielonia Białystok"
$ws.Range("A32").Value = "Legia Warszawa"
$ws.Range("A35").Value = "Remis"
$ws.Range("A39").Value = "Warta Poznań"
$ws.Range("A40").Value = "Remis"
$ws.Range("A41").Value = "Miedź Legnica"
$ws.Range("A42").Value = "Górnik Zabrze"
$ws.Range("A47").Value = "Stal Mielec"
$ws.Range("A48").Value = "Jagielonia Białystok"
$ws.Range("A58").Value = "Lechia Gdańsk"
$ws.Range("A60").Value = "Pogoń Szczecin"
$ws.Range("A62").Value = "Legia Warszawa"
$ws.Range("A63").Value = "Remis"
$ws.Range("A65").Value = "Remis"
$ws.Range("A66").Value = "Miedź Legnica"
$ws.Range("A68").Value = "Remis"
$ws.Range("A72").Value = "Remis"
$ws.Range("A74").Value = "Jagielonia Białystok"
$ws.Range("A75").Value = "Remis"
$ws.Range("A78").Value = "Remis"
$ws.Range("A79").Value = "Remis"
$ws.Range("A80").Value = "Remis"
$ws.Range("A81").Value = "Remis"
$ws.Range("A82").Value = "Śląsk Wrocław"
$ws.Range("A84").Value = "Piast Gliwice"
$ws.Range("A88").Value = "Miedź Legnica"
$ws.Range("A91").Value = "Lech Poznań"
$ws.Range("A92").Value = "Korona Kielce"
$ws.Range("A93").Value = "Legia Warszawa"
$ws.Range("A95").Value = "Remis"
$ws.Range("A96").Value = "Remis"
$ws.Range("A99").Value = "Wisła Płock"
$ws.Range("A100").Value = "Remis"
$ws.Range("A101").Value = "Remis"
$ws.Range("A103").Value = "Remis"
$ws.Range("A108").Value = "Pogoń Szczecin"
$ws.Range("A109").Value = "Remis"
$ws.Range("A111").Value = "Remis"
$ws.Range("A112").Value = "Cracovia"
$ws.Range("A113").Value = "Remis"
$ws.Range("A115").Value = "Śląsk Wrocław"
$ws.Range("A117").Value = "Remis"
$ws.Range("A118").Value = "Remis"
$ws.Range("A123").Value = "Raków Częstochowa"
$ws.Range("A125").Value = "Remis"
$ws.Range("A127").Value = "Remis"
$ws.Range("A129").Value = "Remis"
$ws.Range("A135").Value = "Remis"
$ws.Range("A137").Value = "Jagielonia Białystok"
$ws.Range("A142").Value = "Pogoń Szczecin"
$ws.Range("A148").Value = "Widzew Łódź"
$ws.Range("A154").Value = "Legia Warszawa"
$ws.Range("A155").Value = "Remis"
$ws.Range("A156").Value = "Jagielonia Białystok"
$ws.Range("A159").Value = "Radomiak Radom"
$ws.Range("A161").Value = "Remis"
$ws.Range("A163").Value = "Remis"
$ws.Range("A164").Value = "Remis"
$ws.Range("A165").Value = "Korona Kielce"
$ws.Range("A167").Value = "Remis"
$ws.Range("A169").Value = "Remis"
$ws.Range("A170").Value = "Raków Częstochowa"
$ws.Range("A171").Value = "Widzew Łódź"
$ws.Range("A174").Value = "Remis"
$ws.Range("A177").Value = "Radomiak Radom"
$ws.Range("A181").Value = "Śląsk Wrocław"
$ws.Range("A182").Value = "Cracovia"
$ws.Range("A184").Value = "Lechia Gdańsk"
$ws.Range("A186").Value = "Miedź Legnica"
$ws.Range("A189").Value = "Remis"
$ws.Range("A190").Value = "Śląsk Wrocław"
$ws.Range("A193").Value = "Zagłębie Lubin"
$ws.Range("A196").Value = "Stal Mielec"
$ws.Range("A200").Value = "Remis"
$ws.Range("A201").Value = "Remis"
$ws.Range("A202").Value = "Lech Poznań"
$ws.Range("A203").Value = "Miedź Legnica"
$ws.Range("A204").Value = "Legia Warszawa"
$ws.Range("A206").Value = "Remis"
$ws.Range("A208").Value = "Warta Poznań"
$ws.Range("A210").Value = "Górnik Zabrze"
$ws.Range("A211").Value = "Pogoń Szczecin"
$ws.Range("A212").Value = "Korona Kielce"
$ws.Range("A216").Value = "Cracovia"
$ws.Range("A219").Value = "Miedź Legnica"
$ws.Range("A222").Value = "Legia Warszawa"
$ws.Range("A224").Value = "Lechia Gdańsk"
$ws.Range("A227").Value = "Widzew Łódź"
$ws.Range("A228").Value = "Górnik Zabrze"
$ws.Range("A229").Value = "Warta Poznań"
$ws.Range("A230").Value = "Korona Kielce"
$ws.Range("A231").Value = "Remis"
$ws.Range("A233").Value = "Raków Częstochowa"
$ws.Range("A237").Value = "Warta Poznań"
$ws.Range("A240").Value = "Cracovia"
$ws.Range("A244").Value = "Piast Gliwice"
$ws.Range("A245").Value = "Remis"
$ws.Range("A246").Value = "Piast Gliwice"
$ws.Range("A247").Value = "Remis"
$ws.Range("A249").Value = "Pogoń Szczecin"
$ws.Range("A250").Value = "Lech Poznań"
$ws.Range("A252").Value = "Remis"
$ws.Range("A254").Value = "Remis"
$ws.Range("A255").Value = "Lechia Gdańsk"
$ws.Range("A257").Value = "Górnik Zabrze"
$ws.Range("A260").Value = "Remis"
$ws.Range("A262").Value = "Wisła Płock"
$ws.Range("A263").Value = "Miedź Legnica"
$ws.Range("A265").Value = "Jagielonia Białystok"
$ws.Range("A268").Value = "Remis"
$ws.Range("A271").Value = "Radomiak Radom"
$ws.Range("A275").Value = "Zagłębie Lubin"
$ws.Range("A277").Value = "Górnik Zabrze"
$ws.Range("A278").Value = "Pogoń Szczecin"
$ws.Range("A279").Value = "Remis"
$ws.Range("A280").Value = "Wisła Płock"
$ws.Range("A282").Value = "Remis"
$ws.Range("A283").Value = "Jagielonia Białystok"
$ws.Range("A285").Value = "Remis"
$ws.Range("A288").Value = "Remis"
$ws.Range("A289").Value = "Remis"
$ws.Range("A290").Value = "Cracovia"
$ws.Range("A291").Value = "Stal Mielec"
$ws.Range("A293").Value = "Lechia Gdańsk"
$ws.Range("A294").Value = "Pogoń Szczecin"
$ws.Range("A295").Value = "Radomiak Radom"
$ws.Range("A299").Value = "Remis"
$ws.Range("A302").Value = "Remis"
$ws.Range("A303").Value = "Miedź Legnica"
$ws.Range("A304").Value = "Remis"
$ws.Range("A306").Value = "Warta Poznań"
$ws.Range("A307").Value = "Remis"
